$wb = $excel.ActiveWorkbook

# Sheet ALC, row 32
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(32, 8).Value = 1796.4
$ws.Cells.Item(32, 10).Value = 1796.4
$ws.Cells.Item(32, 12).Value = 1796.4
$ws.Cells.Item(32, 14).Value = -2448.4

# Sheet ALC, row 41
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(41, 8).Value = 2080.9092
$ws.Cells.Item(41, 9).Value = 1670
$ws.Cells.Item(41, 10).Value = 2800
$ws.Cells.Item(41, 11).Value = 1670
$ws.Cells.Item(41, 12).Value = 2800
$ws.Cells.Item(41, 13).Value = -1230
$ws.Cells.Item(41, 14).Value = -3680

# Sheet ALC, row 43
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(43, 8).Value = 11111611
$ws.Cells.Item(43, 9).Value = 499.5
$ws.Cells.Item(43, 11).Value = 499.5
$ws.Cells.Item(43, 13).Value = -430.5

# Sheet ALC, row 113
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(113, 8).Value = 3062.182
$ws.Cells.Item(113, 9).Value = 2526.2856
$ws.Cells.Item(113, 10).Value = 4000
$ws.Cells.Item(113, 11).Value = 2526.2856
$ws.Cells.Item(113, 12).Value = 4000
$ws.Cells.Item(113, 13).Value = 727.7143999999998
$ws.Cells.Item(113, 14).Value = -10508

# Sheet ALC, row 116
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(116, 8).Value = 2937.0588
$ws.Cells.Item(116, 9).Value = 2815.8262
$ws.Cells.Item(116, 11).Value = 2815.8262
$ws.Cells.Item(116, 13).Value = 626.1738

# Sheet ALC, row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(137, 8).Value = 1376.0358
$ws.Cells.Item(137, 9).Value = 997
$ws.Cells.Item(137, 10).Value = 2176.2222
$ws.Cells.Item(137, 11).Value = 2991
$ws.Cells.Item(137, 12).Value = 6528.6666
$ws.Cells.Item(137, 13).Value = -441
$ws.Cells.Item(137, 14).Value = -11628.6666

# Sheet ALC, row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(138, 8).Value = 1461.293
$ws.Cells.Item(138, 10).Value = 1883.209
$ws.Cells.Item(138, 12).Value = 5649.627
$ws.Cells.Item(138, 14).Value = -15929.627

# Sheet ARM, row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 4196.34
$ws.Cells.Item(32, 9).Value = 4220.024
$ws.Cells.Item(32, 11).Value = 4220.024
$ws.Cells.Item(32, 13).Value = -3933.024

# Sheet ARM, row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 1346.6875
$ws.Cells.Item(61, 9).Value = 854.9
$ws.Cells.Item(61, 11).Value = 854.9
$ws.Cells.Item(61, 13).Value = -642.9

# Sheet ARM, row 76
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(76, 8).Value = 29999
$ws.Cells.Item(76, 10).Value = 29999
$ws.Cells.Item(76, 12).Value = 29999
$ws.Cells.Item(76, 14).Value = -30675

# Sheet ARM, row 79
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(79, 8).Value = 29999
$ws.Cells.Item(79, 10).Value = 29999
$ws.Cells.Item(79, 12).Value = 29999
$ws.Cells.Item(79, 14).Value = -32339

# Sheet ARM, row 97
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(97, 8).Value = 531.8182
$ws.Cells.Item(97, 9).Value = 505.55554
$ws.Cells.Item(97, 10).Value = 650
$ws.Cells.Item(97, 11).Value = 505.55554
$ws.Cells.Item(97, 12).Value = 650
$ws.Cells.Item(97, 13).Value = -9.555540000000008
$ws.Cells.Item(97, 14).Value = -1642

# Sheet ARM, row 112
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(112, 8).Value = 29800
$ws.Cells.Item(112, 10).Value = 29800
$ws.Cells.Item(112, 12).Value = 29800
$ws.Cells.Item(112, 14).Value = -32754

# Sheet ARM, row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(136, 8).Value = 1346.6875
$ws.Cells.Item(136, 9).Value = 854.9
$ws.Cells.Item(136, 11).Value = 2564.7
$ws.Cells.Item(136, 13).Value = -14.69999999999982

# Sheet BSM, row 99
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 31251264
$ws.Cells.Item(99, 9).Value = 41667796
$ws.Cells.Item(99, 11).Value = 41667796
$ws.Cells.Item(99, 13).Value = -41666298

# Sheet BSM, row 112
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(112, 8).Value = 34568
$ws.Cells.Item(112, 10).Value = 34568
$ws.Cells.Item(112, 12).Value = 34568
$ws.Cells.Item(112, 14).Value = -37522

# Sheet CRP, row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1274.6812
$ws.Cells.Item(31, 9).Value = 1120.4166
$ws.Cells.Item(31, 10).Value = 2303.111
$ws.Cells.Item(31, 11).Value = 1120.4166
$ws.Cells.Item(31, 12).Value = 2303.111
$ws.Cells.Item(31, 13).Value = -825.4166
$ws.Cells.Item(31, 14).Value = -2893.111

# Sheet CRP, row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(34, 8).Value = 1274.6812
$ws.Cells.Item(34, 9).Value = 1120.4166
$ws.Cells.Item(34, 10).Value = 2303.111
$ws.Cells.Item(34, 11).Value = 1120.4166
$ws.Cells.Item(34, 12).Value = 2303.111
$ws.Cells.Item(34, 13).Value = -918.4166
$ws.Cells.Item(34, 14).Value = -2707.111

# Sheet CRP, row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 1266.3158
$ws.Cells.Item(58, 10).Value = 2250
$ws.Cells.Item(58, 12).Value = 2250
$ws.Cells.Item(58, 14).Value = -2656

# Sheet CRP, row 74
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(74, 8).Value = 32999
$ws.Cells.Item(74, 10).Value = 32999
$ws.Cells.Item(74, 12).Value = 32999
$ws.Cells.Item(74, 14).Value = -34747

# Sheet CRP, row 77
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(77, 8).Value = 32999
$ws.Cells.Item(77, 10).Value = 32999
$ws.Cells.Item(77, 12).Value = 98997
$ws.Cells.Item(77, 14).Value = -107733

# Sheet CRP, row 81
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(81, 8).Value = 16500
$ws.Cells.Item(81, 10).Value = 16500
$ws.Cells.Item(81, 12).Value = 16500
$ws.Cells.Item(81, 14).Value = -18496

# Sheet CRP, row 84
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(84, 8).Value = 16500
$ws.Cells.Item(84, 10).Value = 16500
$ws.Cells.Item(84, 12).Value = 49500
$ws.Cells.Item(84, 14).Value = -59484

# Sheet CRP, row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(132, 8).Value = 2115.75
$ws.Cells.Item(132, 9).Value = 1501.5264
$ws.Cells.Item(132, 10).Value = 3412.4443
$ws.Cells.Item(132, 11).Value = 4504.5792
$ws.Cells.Item(132, 12).Value = 10237.3329
$ws.Cells.Item(132, 13).Value = -1974.5792
$ws.Cells.Item(132, 14).Value = -15297.3329

# Sheet CRP, row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(136, 8).Value = 1266.3158
$ws.Cells.Item(136, 10).Value = 2250
$ws.Cells.Item(136, 12).Value = 6750
$ws.Cells.Item(136, 14).Value = -11850

# Sheet CUL, row 103
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(103, 8).Value = 1427.6154
$ws.Cells.Item(103, 9).Value = 262.5
$ws.Cells.Item(103, 10).Value = 1945.4445
$ws.Cells.Item(103, 11).Value = 787.5
$ws.Cells.Item(103, 12).Value = 5836.333500000001
$ws.Cells.Item(103, 13).Value = 91.5
$ws.Cells.Item(103, 14).Value = -7594.333500000001

# Sheet CUL, row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 15387139
$ws.Cells.Item(131, 10).Value = 2753.39
$ws.Cells.Item(131, 12).Value = 8260.17
$ws.Cells.Item(131, 14).Value = -18340.17

# Sheet CUL, row 139
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(139, 8).Value = 1693.5
$ws.Cells.Item(139, 9).Value = 1886.45
$ws.Cells.Item(139, 10).Value = 1479.1111
$ws.Cells.Item(139, 11).Value = 5659.35
$ws.Cells.Item(139, 12).Value = 4437.3333
$ws.Cells.Item(139, 13).Value = -519.3500000000004
$ws.Cells.Item(139, 14).Value = -14717.3333

# Sheet CUL, row 140
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(140, 8).Value = 33114.63
$ws.Cells.Item(140, 9).Value = 45171.2
$ws.Cells.Item(140, 10).Value = 2973.2
$ws.Cells.Item(140, 11).Value = 135513.6
$ws.Cells.Item(140, 12).Value = 8919.599999999999
$ws.Cells.Item(140, 13).Value = -130333.6
$ws.Cells.Item(140, 14).Value = -19279.6

# Sheet LTW, row 101
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(101, 8).Value = 16998.666
$ws.Cells.Item(101, 10).Value = 16998.666
$ws.Cells.Item(101, 12).Value = 16998.666
$ws.Cells.Item(101, 14).Value = -23488.666

# Sheet WVR, row 68
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(68, 8).Value = 61000
$ws.Cells.Item(68, 10).Value = 61000
$ws.Cells.Item(68, 12).Value = 61000
$ws.Cells.Item(68, 14).Value = -62622

# Sheet WVR, row 71
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(71, 8).Value = 61000
$ws.Cells.Item(71, 10).Value = 61000
$ws.Cells.Item(71, 12).Value = 183000
$ws.Cells.Item(71, 14).Value = -191112

# Sheet WVR, row 119
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(119, 8).Value = 22065.334
$ws.Cells.Item(119, 10).Value = 22065.334
$ws.Cells.Item(119, 12).Value = 22065.334
$ws.Cells.Item(119, 14).Value = -31741.334
